# Insert a new data row at row 2 (pushing the existing rows down by one),
# then populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Insert()

# Reset formatting on the freshly inserted row to match the plain (unstyled)
# data rows, then re-apply the date number format to the Fecha column (D).
$ws.Range("A2:R2").Style = "Normal"
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(2, 3).Value = "Ñuble"
$ws.Cells.Item(2, 4).Value = 44860
$ws.Cells.Item(2, 5).Value = 16
$ws.Cells.Item(2, 6).Value = 300000000
$ws.Cells.Item(2, 7).Value = "Espárragos"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 1200
$ws.Cells.Item(2, 11).Value = 1000
$ws.Cells.Item(2, 12).Value = 1200
$ws.Cells.Item(2, 13).Value = 1100
$ws.Cells.Item(2, 14).Value = "$/kilo"
$ws.Cells.Item(2, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(2, 16).Value = 1100
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = "Hortaliza"
